$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H19").Value = 1320.1428
$ws.Range("J19").Value = 1281.6666
$ws.Range("L19").Value = 1281.6666
$ws.Range("N19").Value = -1631.6666
$ws.Range("H33").Value = 4998.2383
$ws.Range("I33").Value = 6021.4707
$ws.Range("J33").Value = 649.5
$ws.Range("K33").Value = 6021.4707
$ws.Range("L33").Value = 649.5
$ws.Range("M33").Value = -5792.4707
$ws.Range("N33").Value = -1107.5
$ws.Range("H34").Value = 7977.8
$ws.Range("I34").Value = 7977.8
$ws.Range("K34").Value = 7977.8
$ws.Range("M34").Value = -7774.8
$ws.Range("H36").Value = 7977.8
$ws.Range("I36").Value = 7977.8
$ws.Range("K36").Value = 7977.8
$ws.Range("M36").Value = -7262.8
$ws.Range("H40").Value = 3594.3333
$ws.Range("I40").Value = 4744.5
$ws.Range("J40").Value = 3489.7727
$ws.Range("K40").Value = 4744.5
$ws.Range("L40").Value = 3489.7727
$ws.Range("M40").Value = -4569.5
$ws.Range("N40").Value = -3839.7727
$ws.Range("H54").Value = 16380
$ws.Range("I54").Value = 16380
$ws.Range("K54").Value = 16380
$ws.Range("M54").Value = -15894
$ws.Range("H80").Value = 1247.4546
$ws.Range("I80").Value = 1035.1666
$ws.Range("J80").Value = 1502.2
$ws.Range("K80").Value = 3105.4998
$ws.Range("L80").Value = 4506.6
$ws.Range("M80").Value = -2107.4998
$ws.Range("N80").Value = -6502.6
$ws.Range("H83").Value = 1247.4546
$ws.Range("I83").Value = 1035.1666
$ws.Range("J83").Value = 1502.2
$ws.Range("K83").Value = 9316.499400000001
$ws.Range("L83").Value = 13519.8
$ws.Range("M83").Value = -4324.499400000001
$ws.Range("N83").Value = -23503.8
$ws.Range("H86").Value = 2057.25
$ws.Range("I86").Value = 2139
$ws.Range("J86").Value = 1975.5
$ws.Range("K86").Value = 2139
$ws.Range("L86").Value = 1975.5
$ws.Range("M86").Value = -1016
$ws.Range("N86").Value = -4221.5
$ws.Range("H89").Value = 2057.25
$ws.Range("I89").Value = 2139
$ws.Range("J89").Value = 1975.5
$ws.Range("K89").Value = 10695
$ws.Range("L89").Value = 9877.5
$ws.Range("M89").Value = -5079
$ws.Range("N89").Value = -21109.5
$ws.Range("H97").Value = 3497.25
$ws.Range("J97").Value = 4396.3335
$ws.Range("L97").Value = 13189.0005
$ws.Range("N97").Value = -14181.0005
$ws.Range("H98").Value = 4557.6
$ws.Range("I98").Value = 3894.5
$ws.Range("J98").Value = 4999.6665
$ws.Range("K98").Value = 3894.5
$ws.Range("L98").Value = 4999.6665
$ws.Range("M98").Value = -2396.5
$ws.Range("N98").Value = -7995.6665
$ws.Range("H100").Value = 6024.4287
$ws.Range("I100").Value = 6318.25
$ws.Range("J100").Value = 5632.6665
$ws.Range("K100").Value = 6318.25
$ws.Range("L100").Value = 5632.6665
$ws.Range("M100").Value = -5777.25
$ws.Range("N100").Value = -6714.6665
$ws.Range("H103").Value = 1483.5
$ws.Range("I103").Value = 1587.875
$ws.Range("J103").Value = 1274.75
$ws.Range("K103").Value = 4763.625
$ws.Range("L103").Value = 3824.25
$ws.Range("M103").Value = -4177.625
$ws.Range("N103").Value = -4996.25
$ws.Range("H106").Value = 4691.421
$ws.Range("I106").Value = 2974.2354
$ws.Range("K106").Value = 2974.2354
$ws.Range("M106").Value = -2343.2354
$ws.Range("H107").Value = 897.2
$ws.Range("I107").Value = 692.6667
$ws.Range("J107").Value = 1715.3334
$ws.Range("K107").Value = 692.6667
$ws.Range("L107").Value = 1715.3334
$ws.Range("M107").Value = 1227.3333
$ws.Range("N107").Value = -5555.3334
$ws.Range("H112").Value = 2841.2307
$ws.Range("J112").Value = 3215.4546
$ws.Range("L112").Value = 9646.363799999999
$ws.Range("N112").Value = -11862.3638
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442
$ws.Range("H122").Value = 4557.6
$ws.Range("I122").Value = 3894.5
$ws.Range("J122").Value = 4999.6665
$ws.Range("K122").Value = 11683.5
$ws.Range("L122").Value = 14998.9995
$ws.Range("M122").Value = -9233.5
$ws.Range("N122").Value = -19898.9995
$ws.Range("H129").Value = 1054.3871
$ws.Range("I129").Value = 698.7273
$ws.Range("K129").Value = 2096.1819
$ws.Range("M129").Value = 2903.8181
$ws.Range("H131").Value = 6197.5713
$ws.Range("I131").Value = 3897.1667
$ws.Range("K131").Value = 11691.5001
$ws.Range("M131").Value = -6651.500100000001
$ws.Range("H132").Value = 16901.334
$ws.Range("I132").Value = 16131.412
$ws.Range("K132").Value = 48394.236
$ws.Range("M132").Value = -45864.236
$ws.Range("H137").Value = 41597.383
$ws.Range("I137").Value = 55349.895
$ws.Range("J137").Value = 4269.143
$ws.Range("K137").Value = 166049.685
$ws.Range("L137").Value = 12807.429
$ws.Range("M137").Value = -163499.685
$ws.Range("N137").Value = -17907.429
$ws.Range("H138").Value = 2821.3125
$ws.Range("J138").Value = 2774.5789
$ws.Range("L138").Value = 8323.736699999999
$ws.Range("N138").Value = -18603.7367
$ws.Range("H141").Value = 1819.0769
$ws.Range("I141").Value = 1624.6666
$ws.Range("K141").Value = 4873.9998
$ws.Range("M141").Value = 306.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3639.625
$ws.Range("I45").Value = 1442.2858
$ws.Range("J45").Value = 4254.88
$ws.Range("K45").Value = 1442.2858
$ws.Range("L45").Value = 4254.88
$ws.Range("M45").Value = -1065.2858
$ws.Range("N45").Value = -5008.88
$ws.Range("H61").Value = 5182.6924
$ws.Range("I61").Value = 5182.6924
$ws.Range("K61").Value = 5182.6924
$ws.Range("M61").Value = -4970.6924
$ws.Range("H74").Value = 1693.15
$ws.Range("I74").Value = 1609.7059
$ws.Range("J74").Value = 2166
$ws.Range("K74").Value = 1609.7059
$ws.Range("L74").Value = 2166
$ws.Range("M74").Value = -735.7058999999999
$ws.Range("N74").Value = -3914
$ws.Range("H77").Value = 1693.15
$ws.Range("I77").Value = 1609.7059
$ws.Range("J77").Value = 2166
$ws.Range("K77").Value = 8048.5295
$ws.Range("L77").Value = 10830
$ws.Range("M77").Value = -3680.5295
$ws.Range("N77").Value = -19566
$ws.Range("H102").Value = 1812.9412
$ws.Range("I102").Value = 1735.8
$ws.Range("J102").Value = 2391.5
$ws.Range("K102").Value = 1735.8
$ws.Range("L102").Value = 2391.5
$ws.Range("M102").Value = -113.8
$ws.Range("N102").Value = -5635.5
$ws.Range("H122").Value = 1404.7273
$ws.Range("I122").Value = 1411.2188
$ws.Range("J122").Value = 1197
$ws.Range("K122").Value = 4233.6564
$ws.Range("L122").Value = 3591
$ws.Range("M122").Value = -1783.6564
$ws.Range("N122").Value = -8491
$ws.Range("H132").Value = 18870.34
$ws.Range("I132").Value = 20351
$ws.Range("K132").Value = 61053
$ws.Range("M132").Value = -58523
$ws.Range("H136").Value = 5182.6924
$ws.Range("I136").Value = 5182.6924
$ws.Range("K136").Value = 15548.0772
$ws.Range("M136").Value = -12998.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 63518
$ws.Range("I22").Value = 77942.84
$ws.Range("J22").Value = 1010.3333
$ws.Range("K22").Value = 77942.84
$ws.Range("L22").Value = 1010.3333
$ws.Range("M22").Value = -77769.84
$ws.Range("N22").Value = -1356.3333
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H45").Value = 32532.5
$ws.Range("J45").Value = 32532.5
$ws.Range("L45").Value = 32532.5
$ws.Range("N45").Value = -34148.5
$ws.Range("H74").Value = 59998.668
$ws.Range("J74").Value = 59998.668
$ws.Range("L74").Value = 59998.668
$ws.Range("N74").Value = -61870.668
$ws.Range("H77").Value = 59998.668
$ws.Range("J77").Value = 59998.668
$ws.Range("L77").Value = 179996.004
$ws.Range("N77").Value = -189356.004
$ws.Range("H82").Value = 14649.777
$ws.Range("I82").Value = 7731
$ws.Range("J82").Value = 70000
$ws.Range("K82").Value = 7731
$ws.Range("L82").Value = 70000
$ws.Range("M82").Value = -7348
$ws.Range("N82").Value = -70766
$ws.Range("H85").Value = 14649.777
$ws.Range("I85").Value = 7731
$ws.Range("J85").Value = 70000
$ws.Range("K85").Value = 7731
$ws.Range("L85").Value = 70000
$ws.Range("M85").Value = -6405
$ws.Range("N85").Value = -72652
$ws.Range("H94").Value = 2081.9443
$ws.Range("I94").Value = 1810.4814
$ws.Range("J94").Value = 2896.3333
$ws.Range("K94").Value = 1810.4814
$ws.Range("L94").Value = 2896.3333
$ws.Range("M94").Value = -1359.4814
$ws.Range("N94").Value = -3798.3333
$ws.Range("H105").Value = 226231
$ws.Range("I105").Value = 4181.5
$ws.Range("J105").Value = 670330
$ws.Range("K105").Value = 4181.5
$ws.Range("L105").Value = 670330
$ws.Range("M105").Value = -2434.5
$ws.Range("N105").Value = -673824
$ws.Range("H107").Value = 3437.375
$ws.Range("I107").Value = 2940
$ws.Range("K107").Value = 2940
$ws.Range("M107").Value = -1020
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678
$ws.Range("H115").Value = 80000
$ws.Range("J115").Value = 80000
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -83134
$ws.Range("H134").Value = 2941.1052
$ws.Range("I134").Value = 2941.1052
$ws.Range("K134").Value = 8823.3156
$ws.Range("M134").Value = -6288.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 766.3333
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 599
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 599
$ws.Range("M16").Value = -563
$ws.Range("N16").Value = -1173
$ws.Range("H22").Value = 1579.1666
$ws.Range("I22").Value = 944.6667
$ws.Range("K22").Value = 944.6667
$ws.Range("M22").Value = -594.6667
$ws.Range("H31").Value = 2585.5173
$ws.Range("J31").Value = 6869
$ws.Range("L31").Value = 6869
$ws.Range("N31").Value = -7459
$ws.Range("H34").Value = 2585.5173
$ws.Range("J34").Value = 6869
$ws.Range("L34").Value = 6869
$ws.Range("N34").Value = -7273
$ws.Range("H58").Value = 29959.055
$ws.Range("I58").Value = 35257.965
$ws.Range("J58").Value = 3464.5
$ws.Range("K58").Value = 35257.965
$ws.Range("L58").Value = 3464.5
$ws.Range("M58").Value = -35054.965
$ws.Range("N58").Value = -3870.5
$ws.Range("H62").Value = 3528.8333
$ws.Range("I62").Value = 3437.6
$ws.Range("J62").Value = 3985
$ws.Range("K62").Value = 3437.6
$ws.Range("L62").Value = 3985
$ws.Range("M62").Value = -2813.6
$ws.Range("N62").Value = -5233
$ws.Range("H65").Value = 3528.8333
$ws.Range("I65").Value = 3437.6
$ws.Range("J65").Value = 3985
$ws.Range("K65").Value = 17188
$ws.Range("L65").Value = 19925
$ws.Range("M65").Value = -14068
$ws.Range("N65").Value = -26165
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 766.3333
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 599
$ws.Range("K113").Value = 850
$ws.Range("L113").Value = 599
$ws.Range("M113").Value = 1320
$ws.Range("N113").Value = -4939
$ws.Range("H132").Value = 2452.238
$ws.Range("I132").Value = 1908.0526
$ws.Range("K132").Value = 5724.1578
$ws.Range("M132").Value = -3194.1578
$ws.Range("H134").Value = 49167.715
$ws.Range("I134").Value = 60561.53
$ws.Range("J134").Value = 744
$ws.Range("K134").Value = 181684.59
$ws.Range("L134").Value = 2232
$ws.Range("M134").Value = -179149.59
$ws.Range("N134").Value = -7302
$ws.Range("H136").Value = 29959.055
$ws.Range("I136").Value = 35257.965
$ws.Range("J136").Value = 3464.5
$ws.Range("K136").Value = 105773.895
$ws.Range("L136").Value = 10393.5
$ws.Range("M136").Value = -103223.895
$ws.Range("N136").Value = -15493.5
$ws.Range("H140").Value = 79865.2
$ws.Range("J140").Value = 79865.2
$ws.Range("L140").Value = 79865.2
$ws.Range("N140").Value = -90225.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 612.6667
$ws.Range("I5").Value = 612.6667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1838.0001
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1726.0001
$ws.Range("N5").ClearContents()
$ws.Range("H11").Value = 110
$ws.Range("I11").Value = 110
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 330
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -190
$ws.Range("N11").ClearContents()
$ws.Range("H38").Value = 884.2
$ws.Range("I38").Value = 390.33334
$ws.Range("J38").Value = 1625
$ws.Range("K38").Value = 1171.00002
$ws.Range("L38").Value = 4875
$ws.Range("M38").Value = -824.0000199999999
$ws.Range("N38").Value = -5569
$ws.Range("H104").Value = 1000
$ws.Range("I104").Value = 1000
$ws.Range("K104").Value = 3000
$ws.Range("M104").Value = -379
$ws.Range("H121").Value = 468.14285
$ws.Range("J121").Value = 514.5
$ws.Range("L121").Value = 1543.5
$ws.Range("N121").Value = -4163.5
$ws.Range("H135").Value = 612.6667
$ws.Range("I135").Value = 612.6667
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5514.0003
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2979.0003
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2988.611
$ws.Range("J80").Value = 3192.5557
$ws.Range("L80").Value = 3192.5557
$ws.Range("N80").Value = -5188.5557
$ws.Range("H83").Value = 2988.611
$ws.Range("J83").Value = 3192.5557
$ws.Range("L83").Value = 15962.7785
$ws.Range("N83").Value = -25946.7785
$ws.Range("H102").Value = 2598.7188
$ws.Range("I102").Value = 2448.8635
$ws.Range("K102").Value = 2448.8635
$ws.Range("M102").Value = -826.8634999999999
$ws.Range("H122").Value = 2214.9285
$ws.Range("I122").Value = 1642.7188
$ws.Range("K122").Value = 4928.1564
$ws.Range("M122").Value = -2478.1564
$ws.Range("H126").Value = 6931.64
$ws.Range("I126").Value = 6743
$ws.Range("J126").Value = 7214.6
$ws.Range("K126").Value = 20229
$ws.Range("L126").Value = 21643.8
$ws.Range("M126").Value = -17759
$ws.Range("N126").Value = -26583.8
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 28996.447
$ws.Range("I132").Value = 44253.918
$ws.Range("K132").Value = 132761.754
$ws.Range("M132").Value = -130231.754

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2482.077
$ws.Range("I7").Value = 2355.5833
$ws.Range("K7").Value = 2355.5833
$ws.Range("M7").Value = -2243.5833
$ws.Range("H22").Value = 80806.5
$ws.Range("J22").Value = 1999.1111
$ws.Range("L22").Value = 1999.1111
$ws.Range("N22").Value = -2589.1111
$ws.Range("H27").Value = 80806.5
$ws.Range("J27").Value = 1999.1111
$ws.Range("L27").Value = 1999.1111
$ws.Range("N27").Value = -2213.1111
$ws.Range("H42").Value = 32666.334
$ws.Range("J42").Value = 38999.5
$ws.Range("L42").Value = 38999.5
$ws.Range("N42").Value = -40125.5
$ws.Range("H46").Value = 34005.43
$ws.Range("I46").Value = 56475
$ws.Range("J46").Value = 4046
$ws.Range("K46").Value = 56475
$ws.Range("L46").Value = 4046
$ws.Range("M46").Value = -56287
$ws.Range("N46").Value = -4422
$ws.Range("H49").Value = 32666.334
$ws.Range("J49").Value = 38999.5
$ws.Range("L49").Value = 38999.5
$ws.Range("N49").Value = -39293.5
$ws.Range("H68").Value = 5497.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5497.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5497.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6995.5
$ws.Range("H71").Value = 5497.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5497.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 27487.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -34975.5
$ws.Range("H80").Value = 69998
$ws.Range("J80").Value = 69998
$ws.Range("L80").Value = 69998
$ws.Range("N80").Value = -72244
$ws.Range("H82").Value = 2965.7144
$ws.Range("I82").Value = 1114.2
$ws.Range("J82").Value = 3994.3333
$ws.Range("K82").Value = 1114.2
$ws.Range("L82").Value = 3994.3333
$ws.Range("M82").Value = -753.2
$ws.Range("N82").Value = -4716.3333
$ws.Range("H83").Value = 69998
$ws.Range("J83").Value = 69998
$ws.Range("L83").Value = 209994
$ws.Range("N83").Value = -221226
$ws.Range("H85").Value = 2965.7144
$ws.Range("I85").Value = 1114.2
$ws.Range("J85").Value = 3994.3333
$ws.Range("K85").Value = 1114.2
$ws.Range("L85").Value = 3994.3333
$ws.Range("M85").Value = 133.8
$ws.Range("N85").Value = -6490.3333
$ws.Range("H100").Value = 3680.8333
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 2987.7273
$ws.Range("I122").Value = 2609.625
$ws.Range("K122").Value = 7828.875
$ws.Range("M122").Value = -5378.875
$ws.Range("H126").Value = 2482.077
$ws.Range("I126").Value = 2355.5833
$ws.Range("K126").Value = 7066.749899999999
$ws.Range("M126").Value = -4596.749899999999
$ws.Range("H132").Value = 64134.35
$ws.Range("I132").Value = 83643.60000000001
$ws.Range("K132").Value = 250930.8
$ws.Range("M132").Value = -248400.8
$ws.Range("H136").Value = 4529.1113
$ws.Range("I136").Value = 3937.8572
$ws.Range("J136").Value = 6598.5
$ws.Range("K136").Value = 11813.5716
$ws.Range("L136").Value = 19795.5
$ws.Range("M136").Value = -9263.571599999999
$ws.Range("N136").Value = -24895.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 179829.83
$ws.Range("J62").Value = 179829.83
$ws.Range("L62").Value = 179829.83
$ws.Range("N62").Value = -181077.83
$ws.Range("H65").Value = 179829.83
$ws.Range("J65").Value = 179829.83
$ws.Range("L65").Value = 899149.1499999999
$ws.Range("N65").Value = -905389.1499999999
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H113").Value = 1370.5714
$ws.Range("I113").Value = 634.125
$ws.Range("K113").Value = 1902.375
$ws.Range("M113").Value = 267.625
$ws.Range("H122").Value = 1726.1111
$ws.Range("I122").Value = 1466.3939
$ws.Range("J122").Value = 4583
$ws.Range("K122").Value = 4399.1817
$ws.Range("L122").Value = 13749
$ws.Range("M122").Value = -1949.1817
$ws.Range("N122").Value = -18649
$ws.Range("H126").Value = 173362.67
$ws.Range("I126").Value = 338131.34
$ws.Range("J126").Value = 8594
$ws.Range("K126").Value = 1014394.02
$ws.Range("L126").Value = 25782
$ws.Range("M126").Value = -1011924.02
$ws.Range("N126").Value = -30722
$ws.Range("H132").Value = 25951.023
$ws.Range("I132").Value = 29498.135
$ws.Range("K132").Value = 88494.405
$ws.Range("M132").Value = -85964.405
$ws.Range("H136").Value = 3108.7693
$ws.Range("I136").Value = 3108.7693
$ws.Range("K136").Value = 9326.3079
$ws.Range("M136").Value = -6776.3079
